$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Cells.Item(1,1).Value2
Write-Output $ws.Cells.Item(2,1).Value2
Write-Output $ws.Range("A1").Value2
Write-Output $ws.Range("A2").Value2
